{"js": "// Remove the \"IT Support Intern\" text that was filled in under the\n// \"Department Assigned:\" table cell, leaving the paragraph mark (and its\n// formatting) intact but with no runs.\nconst results = context.document.body.search(\"IT Support Intern\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Remove the leftover \"_GoBack\" bookmark (start + end) left over from the\n// author's last cursor position - Word normally strips this on save.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the \"IT Support Intern\" text that was filled in under the\n# \"Department Assigned:\" table cell, leaving the paragraph mark (and its\n# formatting) intact but with no runs.\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"IT Support Intern\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\nif ($find.Found) {\n    $range.Text = \"\"\n}\n\n# Remove the leftover \"_GoBack\" bookmark left over from the author's last\n# cursor position - Word normally strips this on save.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
